$d = $word.ActiveDocument

$replacements = @(
    @{old="446×2=892"; new="675×8=5400"},
    @{old="780×7=5460"; new="702×5=3510"},
    @{old="133×9=1197"; new="947×7=6629"},
    @{old="799×7=5593"; new="195×2=390"},
    @{old="641×3=1923"; new="831×3=2493"},
    @{old="298×4=1192"; new="496×6=2976"},
    @{old="524×7=3668"; new="438×5=2190"},
    @{old="364×3=1092"; new="523×8=4184"},
    @{old="304×7=2128"; new="935×8=7480"},
    @{old="592×4=2368"; new="510×8=4080"},
    @{old="791×4=3164"; new="468×2=936"},
    @{old="945×5=4725"; new="356×6=2136"},
    @{old="348×8=2784"; new="226×7=1582"},
    @{old="994×9=8946"; new="420×5=2100"},
    @{old="587×7=4109"; new="822×6=4932"},
    @{old="949×2=1898"; new="316×9=2844"},
    @{old="855×3=2565"; new="462×7=3234"},
    @{old="310×8=2480"; new="944×8=7552"},
    @{old="825×5=4125"; new="391×3=1173"},
    @{old="804×6=4824"; new="891×4=3564"},
    @{old="892×8=7136"; new="754×3=2262"},
    @{old="948×8=7584"; new="528×2=1056"},
    @{old="996×5=4980"; new="301×5=1505"},
    @{old="103×7=721"; new="235×3=705"},
    @{old="945×7=6615"; new="259×4=1036"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
